$wb = $excel.ActiveWorkbook

# --- ALC sheet: row 137 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1023.81635
$ws.Range("I137").Value = 914.125
$ws.Range("J137").Value = 1511.3334
$ws.Range("K137").Value = 2742.375
$ws.Range("L137").Value = 4534.0002
$ws.Range("M137").Value = -192.375
$ws.Range("N137").Value = -9634.0002

# --- CRP sheet: rows 31, 34, 99, 126 ---
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 1239.0968
$ws.Range("I31").Value = 847
$ws.Range("J31").Value = 1541.5714
$ws.Range("K31").Value = 847
$ws.Range("L31").Value = 1541.5714
$ws.Range("M31").Value = -552
$ws.Range("N31").Value = -2131.5714

$ws.Range("H34").Value = 1239.0968
$ws.Range("I34").Value = 847
$ws.Range("J34").Value = 1541.5714
$ws.Range("K34").Value = 847
$ws.Range("L34").Value = 1541.5714
$ws.Range("M34").Value = -645
$ws.Range("N34").Value = -1945.5714

$ws.Range("H99").Value = 76925990
$ws.Range("I99").Value = 142859310
$ws.Range("J99").Value = 3783.3333
$ws.Range("K99").Value = 142859310
$ws.Range("L99").Value = 3783.3333
$ws.Range("M99").Value = -142857812
$ws.Range("N99").Value = -6779.3333

$ws.Range("H126").Value = 76925990
$ws.Range("I126").Value = 142859310
$ws.Range("J126").Value = 3783.3333
$ws.Range("K126").Value = 428577930
$ws.Range("L126").Value = 11349.9999
$ws.Range("M126").Value = -428575460
$ws.Range("N126").Value = -16289.9999

# --- WVR sheet: rows 119-141 (except 134) lose their H:N values entirely ---
$ws = $wb.Worksheets.Item("WVR")
$rowsToClear = @(119,120,121,122,123,124,125,126,127,128,129,130,131,132,133,135,136,137,138,139,140,141)
foreach ($r in $rowsToClear) {
    $ws.Range("H" + $r + ":N" + $r).ClearContents()
}
